$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.888.56'
$ws.Range('E2').Value = '  +1.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.313.52'
$ws.Range('E3').Value = '  +1.35%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '187.70'
$ws.Range('E5').Value = '  +5.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '554.51'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.310.70'
$ws.Range('E8').Value = '  +1.39%  '
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('E10').Value = '  -3.83%  '
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000267'
$ws.Range('E13').Value = '  +2.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.60'
$ws.Range('E14').Value = '  +1.32%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.841.44'
$ws.Range('E15').Value = '  +1.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '596.88'
$ws.Range('E16').Value = '  -1.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.890.10'
$ws.Range('E17').Value = '  +1.15%  '
$ws.Range('E18').Value = '  +1.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.89'
$ws.Range('E19').Value = '  +0.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.325.58'
$ws.Range('E20').Value = '  +1.88%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.99'
$ws.Range('E21').Value = '  -2.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.897'
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.35'
$ws.Range('E23').Value = '  +4.99%  '
$ws.Range('E24').Value = '  +1.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '99.88'
$ws.Range('E25').Value = '  -1.37%  '
$ws.Range('E26').Value = '  -0.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.06'
$ws.Range('E27').Value = '  +1.36%  '
$ws.Range('E28').Value = '  +2.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.46'
$ws.Range('E29').Value = '  +1.64%  '
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.36'
$ws.Range('E31').Value = '  +0.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.68'
$ws.Range('E32').Value = '  +7.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.83'
$ws.Range('E33').Value = '  -0.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '577.66'
$ws.Range('E34').Value = '  +7.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '10.99'
$ws.Range('E35').Value = '  +0.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.104'
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.705.17'
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '56.65'
$ws.Range('E39').Value = '  +0.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.51'
$ws.Range('E40').Value = '  +8.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '33.53'
$ws.Range('E41').Value = '  +6.01%  '
$ws.Range('E42').Value = '  +1.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.22'
$ws.Range('E43').Value = '  -5.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0₃0697'
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.43'
$ws.Range('E45').Value = '  +8.06%  '
$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.63'
$ws.Range('E46').Value = '  -0.87%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.336'
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0415'
$ws.Range('E48').Value = '  +2.06%  '
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.55'
$ws.Range('E50').Value = '  -0.64%  '
$ws.Range('B51').Value = 'FirstDigitalUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.998'
$ws.Range('E51').Value = '  +0.06%  '
